$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure numeric-looking price cells stay text (match existing inline-string cells)
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

# Apply the updated values
$ws.Range('D2').Value = '51.028.92'
$ws.Range('E2').Value = '  -0.74%  '
$ws.Range('D3').Value = '2.932.76'
$ws.Range('E3').Value = '  -1.04%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '373.76'
$ws.Range('E5').Value = '  -1.00%  '
$ws.Range('D6').Value = '100.38'
$ws.Range('E6').Value = '  -3.76%  '
$ws.Range('D7').Value = '0.532'
$ws.Range('E7').Value = '  -1.25%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').Value = '0.579'
$ws.Range('E9').Value = '  -1.72%  '
$ws.Range('D10').Value = '36.05'
$ws.Range('E10').Value = '  -2.48%  '
$ws.Range('D11').Value = '0.139'
$ws.Range('E11').Value = '  -0.54%  '
$ws.Range('D12').Value = '0.0845'
$ws.Range('E12').Value = '  +0.60%  '
$ws.Range('D13').Value = '3.396.14'
$ws.Range('E13').Value = '  -1.07%  '
$ws.Range('D14').Value = '17.97'
$ws.Range('E14').Value = '  -1.75%  '
$ws.Range('D15').Value = '7.45'
$ws.Range('E15').Value = '  -1.12%  '
$ws.Range('D16').Value = '11.34'
$ws.Range('E16').Value = '  +53.93%  '
$ws.Range('D17').Value = '2.938.04'
$ws.Range('E17').Value = '  -1.01%  '
$ws.Range('D18').Value = '0.969'
$ws.Range('E18').Value = '  +1.11%  '
$ws.Range('D19').Value = '51.020.49'
$ws.Range('E19').Value = '  -0.66%  '
$ws.Range('E20').Value = '  -5.42%  '
$ws.Range('D21').Value = '12.29'
$ws.Range('E21').Value = '  -4.47%  '
$ws.Range('D22').Value = '0.0₃0952'
$ws.Range('E22').Value = '  -0.64%  '
$ws.Range('D23').Value = '263.66'
$ws.Range('E23').Value = '  +1.03%  '
$ws.Range('D24').Value = '68.43'
$ws.Range('E24').Value = '  -1.17%  '
$ws.Range('D25').Value = '3.10'
$ws.Range('E25').Value = '  +10.10%  '
$ws.Range('D26').Value = '8.09'
$ws.Range('E26').Value = '  -0.38%  '
$ws.Range('D27').Value = '7.33'
$ws.Range('E27').Value = '  -2.30%  '
$ws.Range('E28').Value = '  +0.11%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '25.49'
$ws.Range('E29').Value = '  -1.21%  '
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').Value = '0.163'
$ws.Range('E30').Value = '  -4.12%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '0.110'
$ws.Range('E31').Value = '  -3.08%  '
$ws.Range('D32').Value = '9.88'
$ws.Range('E32').Value = '  +0.43%  '
$ws.Range('D33').Value = '50.55'
$ws.Range('E33').Value = '  -0.79%  '
$ws.Range('D34').Value = '2.04'
$ws.Range('E34').Value = '  -1.92%  '
$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D35').Value = '0.0440'
$ws.Range('E35').Value = '  -0.74%  '
$ws.Range('B36').Value = 'InjectiveProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D36').Value = '32.88'
$ws.Range('E36').Value = '  -5.57%  '
$ws.Range('E37').Value = '  -0.19%  '
$ws.Range('D38').Value = '3.14'
$ws.Range('E38').Value = '  +4.27%  '
$ws.Range('D39').Value = '0.115'
$ws.Range('E39').Value = '  -0.31%  '
$ws.Range('D40').Value = '16.27'
$ws.Range('E40').Value = '  -4.72%  '
$ws.Range('D41').Value = '1.78'
$ws.Range('E41').Value = '  -3.38%  '
$ws.Range('D42').Value = '2.45'
$ws.Range('E42').Value = '  -4.66%  '
$ws.Range('D43').Value = '119.57'
$ws.Range('E43').Value = '  -4.20%  '
$ws.Range('D44').Value = '20.98'
$ws.Range('E44').Value = '  -2.75%  '
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').Value = '0.275'
$ws.Range('E45').Value = '  -3.92%  '
$ws.Range('B46').Value = 'WEMIXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').Value = '2.02'
$ws.Range('E46').Value = '  -1.51%  '
$ws.Range('D47').Value = '3.27'
$ws.Range('E47').Value = '  +2.10%  '
$ws.Range('D48').Value = '2.30'
$ws.Range('E48').Value = '  -2.83%  '
$ws.Range('D49').Value = '1.978.78'
$ws.Range('E49').Value = '  -2.43%  '
$ws.Range('D50').Value = '0.0324'
$ws.Range('E50').Value = '  -4.19%  '
$ws.Range('D51').Value = '1.31'
$ws.Range('E51').Value = '  +2.12%  '

# Restore default (unformatted) style so cells match original formatting
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').Style = "Normal"
$ws.Range('D7').Style = "Normal"
$ws.Range('D9').Style = "Normal"
$ws.Range('D10').Style = "Normal"
$ws.Range('D11').Style = "Normal"
$ws.Range('D12').Style = "Normal"
$ws.Range('D14').Style = "Normal"
$ws.Range('D15').Style = "Normal"
$ws.Range('D16').Style = "Normal"
$ws.Range('D18').Style = "Normal"
$ws.Range('D21').Style = "Normal"
$ws.Range('D23').Style = "Normal"
$ws.Range('D24').Style = "Normal"
$ws.Range('D25').Style = "Normal"
$ws.Range('D26').Style = "Normal"
$ws.Range('D27').Style = "Normal"
$ws.Range('D29').Style = "Normal"
$ws.Range('D30').Style = "Normal"
$ws.Range('D31').Style = "Normal"
$ws.Range('D32').Style = "Normal"
$ws.Range('D33').Style = "Normal"
$ws.Range('D34').Style = "Normal"
$ws.Range('D35').Style = "Normal"
$ws.Range('D36').Style = "Normal"
$ws.Range('D38').Style = "Normal"
$ws.Range('D39').Style = "Normal"
$ws.Range('D40').Style = "Normal"
$ws.Range('D41').Style = "Normal"
$ws.Range('D42').Style = "Normal"
$ws.Range('D43').Style = "Normal"
$ws.Range('D44').Style = "Normal"
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').Style = "Normal"
$ws.Range('D47').Style = "Normal"
$ws.Range('D48').Style = "Normal"
$ws.Range('D50').Style = "Normal"
$ws.Range('D51').Style = "Normal"
